# Översikt UPPLANDS VÄSBY - automatic update of files.
# A new case (A 5779-2026) was logged, which pushes in as the new row 2 and
# shifts every other case down by one row. The "Förändrad" column (C) is
# bumped to the new run date (46064) for every case, and a few cases swap
# their relative order in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 2; this shifts all existing rows (and their
# styles/number formats) down by one, so every case currently in the sheet
# keeps its own formatting.
$ws.Rows.Item(2).Insert()

# --- Row 2: brand new case, A 5779-2026 -------------------------------
$ws.Range("A2").Value = "A 5779-2026"
$ws.Range("B2:C2").NumberFormat = "YYYY-MM-DD"
$ws.Range("B2").Value = "2026-02-03"
$ws.Range("C2").Value = "2026-02-11"
$ws.Range("D2").Value = "STOCKHOLMS LÄN"
$ws.Range("E2").Value = "UPPLANDS VÄSBY"
$ws.Range("G2").Value = 0.8
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = "Tallticka`nBlåsippa"
$ws.Range("R2").WrapText = $true
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0114/artfynd/A 5779-2026 artfynd.xlsx", "A 5779-2026")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0114/kartor/A 5779-2026 karta.png", "A 5779-2026")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0114/klagomål/A 5779-2026 FSC-klagomål.docx", "A 5779-2026")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0114/klagomålsmail/A 5779-2026 FSC-klagomål mail.docx", "A 5779-2026")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0114/tillsyn/A 5779-2026 tillsynsbegäran.docx", "A 5779-2026")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0114/tillsynsmail/A 5779-2026 tillsynsbegäran mail.docx", "A 5779-2026")'

# Keep the original row height (15pt) instead of the auto-fit height Excel
# would otherwise apply because of the two-line wrapped species list.
$ws.Rows.Item(2).RowHeight = 15

# --- Rows 3-14: existing cases (re-ordered), each refreshed in place ---
# Columns D/E (län/kommun) and H-Q (species counters, all zero) are already
# correct after the row-insert shift, so only Beteckning/Datum/Area/Förändrad
# need to be (re)written per row.

$ws.Range("A3").Value = "A 1053-2022"
$ws.Range("B3").Value = "2022-01-10"
$ws.Range("C3").Value = "2026-02-11"
$ws.Range("G3").Value = 1.7

$ws.Range("A4").Value = "A 24221-2021"
$ws.Range("B4").Value2 = 44336.78922453704
$ws.Range("C4").Value = "2026-02-11"
$ws.Range("G4").Value = 1.1

$ws.Range("A5").Value = "A 46826-2022"
$ws.Range("B5").Value = "2022-10-14"
$ws.Range("C5").Value = "2026-02-11"
$ws.Range("G5").Value = 4.5

$ws.Range("A6").Value = "A 50934-2024"
$ws.Range("B6").Value = "2024-11-06"
$ws.Range("C6").Value = "2026-02-11"
$ws.Range("G6").Value = 0.6

$ws.Range("A7").Value = "A 46779-2025"
$ws.Range("B7").Value = "2025-09-26"
$ws.Range("C7").Value = "2026-02-11"
$ws.Range("G7").Value = 1.5

$ws.Range("A8").Value = "A 31120-2023"
$ws.Range("B8").Value = "2023-07-06"
$ws.Range("C8").Value = "2026-02-11"
$ws.Range("G8").Value = 0.2

$ws.Range("A9").Value = "A 64431-2023"
$ws.Range("B9").Value = "2023-12-20"
$ws.Range("C9").Value = "2026-02-11"
$ws.Range("G9").Value = 0.5

$ws.Range("A10").Value = "A 56948-2025"
$ws.Range("B10").Value2 = 45978.64356481482
$ws.Range("C10").Value = "2026-02-11"
$ws.Range("G10").Value = 4.7

$ws.Range("A11").Value = "A 56917-2025"
$ws.Range("B11").Value2 = 45978.58453703704
$ws.Range("C11").Value = "2026-02-11"
$ws.Range("G11").Value = 0.7

$ws.Range("A12").Value = "A 27724-2022"
$ws.Range("B12").Value2 = 44743.48386574074
$ws.Range("C12").Value = "2026-02-11"
$ws.Range("G12").Value = 1.3

$ws.Range("A13").Value = "A 8748-2022"
$ws.Range("B13").Value = "2022-02-21"
$ws.Range("C13").Value = "2026-02-11"
$ws.Range("G13").Value = 1

$ws.Range("A14").Value = "A 64445-2023"
$ws.Range("B14").Value = "2023-12-20"
$ws.Range("C14").Value = "2026-02-11"
$ws.Range("G14").Value = 3.7
